$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1190.3334
$ws.Range("I58").Value = 356
$ws.Range("J58").Value = 2233.25
$ws.Range("K58").Value = 1068
$ws.Range("L58").Value = 6699.75
$ws.Range("M58").Value = -918
$ws.Range("N58").Value = -6999.75

$ws.Range("H62").Value = 6186.375
$ws.Range("I62").Value = 7899
$ws.Range("J62").Value = 3332
$ws.Range("K62").Value = 7899
$ws.Range("L62").Value = 3332
$ws.Range("M62").Value = -7275
$ws.Range("N62").Value = -4580

$ws.Range("H65").Value = 6186.375
$ws.Range("I65").Value = 7899
$ws.Range("J65").Value = 3332
$ws.Range("K65").Value = 39495
$ws.Range("L65").Value = 16660
$ws.Range("M65").Value = -36375
$ws.Range("N65").Value = -22900

$ws.Range("H92").Value = 383.6316
$ws.Range("I92").Value = 417.93332
$ws.Range("K92").Value = 417.93332
$ws.Range("M92").Value = 830.06668

$ws.Range("H96").Value = 62507264
$ws.Range("I96").Value = 4036.2222
$ws.Range("J96").Value = 142868560
$ws.Range("K96").Value = 12108.6666
$ws.Range("L96").Value = 428605680
$ws.Range("M96").Value = -10735.6666
$ws.Range("N96").Value = -428608426

$ws.Range("H132").Value = 19640.96
$ws.Range("I132").Value = 2954.5813
$ws.Range("J132").Value = 99364.78
$ws.Range("K132").Value = 8863.743899999999
$ws.Range("L132").Value = 298094.34
$ws.Range("M132").Value = -6333.743899999999
$ws.Range("N132").Value = -303154.34

$ws.Range("H138").Value = 1585.21
$ws.Range("I138").Value = 662.1163
$ws.Range("J138").Value = 2281.5789
$ws.Range("K138").Value = 1986.3489
$ws.Range("L138").Value = 6844.736699999999
$ws.Range("M138").Value = 3153.6511
$ws.Range("N138").Value = -17124.7367

$ws.Range("H141").Value = 1790
$ws.Range("I141").Value = 762.8108
$ws.Range("J141").Value = 5590.6
$ws.Range("K141").Value = 2288.4324
$ws.Range("L141").Value = 16771.8
$ws.Range("M141").Value = 2891.5676
$ws.Range("N141").Value = -27131.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10197.691
$ws.Range("I32").Value = 9546.018
$ws.Range("J32").Value = 13977.4
$ws.Range("K32").Value = 9546.018
$ws.Range("L32").Value = 13977.4
$ws.Range("M32").Value = -9259.018
$ws.Range("N32").Value = -14551.4

$ws.Range("H61").Value = 1217.7018
$ws.Range("I61").Value = 1035.4445
$ws.Range("J61").Value = 4498.3335
$ws.Range("K61").Value = 1035.4445
$ws.Range("L61").Value = 4498.3335
$ws.Range("M61").Value = -823.4445000000001
$ws.Range("N61").Value = -4922.3335

$ws.Range("H74").Value = 1510.6666
$ws.Range("I74").Value = 1403.122
$ws.Range("J74").Value = 1951.6
$ws.Range("K74").Value = 1403.122
$ws.Range("L74").Value = 1951.6
$ws.Range("M74").Value = -529.1220000000001
$ws.Range("N74").Value = -3699.6

$ws.Range("H77").Value = 1510.6666
$ws.Range("I77").Value = 1403.122
$ws.Range("J77").Value = 1951.6
$ws.Range("K77").Value = 7015.610000000001
$ws.Range("L77").Value = 9758
$ws.Range("M77").Value = -2647.610000000001
$ws.Range("N77").Value = -18494

$ws.Range("H97").Value = 919.9375
$ws.Range("I97").Value = 609
$ws.Range("J97").Value = 1852.75
$ws.Range("K97").Value = 609
$ws.Range("L97").Value = 1852.75
$ws.Range("M97").Value = -113
$ws.Range("N97").Value = -2844.75

$ws.Range("H102").Value = 14683.444
$ws.Range("I102").Value = 2710
$ws.Range("J102").Value = 24262.2
$ws.Range("K102").Value = 2710
$ws.Range("L102").Value = 24262.2
$ws.Range("M102").Value = -1088
$ws.Range("N102").Value = -27506.2

$ws.Range("H132").Value = 11365254
$ws.Range("I132").Value = 13889877
$ws.Range("J132").Value = 4446.75
$ws.Range("K132").Value = 41669631
$ws.Range("L132").Value = 13340.25
$ws.Range("M132").Value = -41667101
$ws.Range("N132").Value = -18400.25

$ws.Range("H136").Value = 1217.7018
$ws.Range("I136").Value = 1035.4445
$ws.Range("J136").Value = 4498.3335
$ws.Range("K136").Value = 3106.3335
$ws.Range("L136").Value = 13495.0005
$ws.Range("M136").Value = -556.3335000000002
$ws.Range("N136").Value = -18595.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1692.5869
$ws.Range("I99").Value = 1555.9048
$ws.Range("J99").Value = 3127.75
$ws.Range("K99").Value = 1555.9048
$ws.Range("L99").Value = 3127.75
$ws.Range("M99").Value = -57.90480000000002
$ws.Range("N99").Value = -6123.75

$ws.Range("H105").Value = 1507.9183
$ws.Range("I105").Value = 1434.3636
$ws.Range("J105").Value = 1567.8518
$ws.Range("K105").Value = 1434.3636
$ws.Range("L105").Value = 1567.8518
$ws.Range("M105").Value = 312.6364000000001
$ws.Range("N105").Value = -5061.8518

$ws.Range("H134").Value = 2527.0989
$ws.Range("I134").Value = 1093.4565
$ws.Range("J134").Value = 3992.6
$ws.Range("K134").Value = 3280.3695
$ws.Range("L134").Value = 11977.8
$ws.Range("M134").Value = -745.3694999999998
$ws.Range("N134").Value = -17047.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2166.63
$ws.Range("I31").Value = 997.9796
$ws.Range("J31").Value = 3289.451
$ws.Range("K31").Value = 997.9796
$ws.Range("L31").Value = 3289.451
$ws.Range("M31").Value = -702.9796
$ws.Range("N31").Value = -3879.451

$ws.Range("H34").Value = 2166.63
$ws.Range("I34").Value = 997.9796
$ws.Range("J34").Value = 3289.451
$ws.Range("K34").Value = 997.9796
$ws.Range("L34").Value = 3289.451
$ws.Range("M34").Value = -795.9796
$ws.Range("N34").Value = -3693.451

$ws.Range("H58").Value = 985.60565
$ws.Range("I58").Value = 725.8412499999999
$ws.Range("J58").Value = 3031.25
$ws.Range("K58").Value = 725.8412499999999
$ws.Range("L58").Value = 3031.25
$ws.Range("M58").Value = -522.8412499999999
$ws.Range("N58").Value = -3437.25

$ws.Range("H122").Value = 201175
$ws.Range("I122").Value = 301075
$ws.Range("K122").Value = 903225
$ws.Range("M122").Value = -900775

$ws.Range("H132").Value = 19031.975
$ws.Range("I132").Value = 911.5441
$ws.Range("J132").Value = 142250.9
$ws.Range("K132").Value = 2734.6323
$ws.Range("L132").Value = 426752.7
$ws.Range("M132").Value = -204.6322999999998
$ws.Range("N132").Value = -431812.7

$ws.Range("H134").Value = 1250.1321
$ws.Range("I134").Value = 968.7659
$ws.Range("J134").Value = 3454.1667
$ws.Range("K134").Value = 2906.2977
$ws.Range("L134").Value = 10362.5001
$ws.Range("M134").Value = -371.2977000000001
$ws.Range("N134").Value = -15432.5001

$ws.Range("H136").Value = 985.60565
$ws.Range("I136").Value = 725.8412499999999
$ws.Range("J136").Value = 3031.25
$ws.Range("K136").Value = 2177.52375
$ws.Range("L136").Value = 9093.75
$ws.Range("M136").Value = 372.4762500000002
$ws.Range("N136").Value = -14193.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3003.0222
$ws.Range("I5").Value = 4134.4814
$ws.Range("J5").Value = 1305.8334
$ws.Range("K5").Value = 12403.4442
$ws.Range("L5").Value = 3917.5002
$ws.Range("M5").Value = -12291.4442
$ws.Range("N5").Value = -4141.5002

$ws.Range("H68").Value = 5630.6665
$ws.Range("J68").Value = 8046
$ws.Range("L68").Value = 24138
$ws.Range("N68").Value = -25760

$ws.Range("H71").Value = 5630.6665
$ws.Range("J71").Value = 8046
$ws.Range("L71").Value = 72414
$ws.Range("N71").Value = -80526

$ws.Range("H113").Value = 5654.2
$ws.Range("J113").Value = 686.1667
$ws.Range("L113").Value = 2058.5001
$ws.Range("N113").Value = -6398.5001

$ws.Range("H135").Value = 3003.0222
$ws.Range("I135").Value = 4134.4814
$ws.Range("J135").Value = 1305.8334
$ws.Range("K135").Value = 37210.33259999999
$ws.Range("L135").Value = 11752.5006
$ws.Range("M135").Value = -34675.33259999999
$ws.Range("N135").Value = -16822.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 6003
$ws.Range("I4").Value = 6003
$ws.Range("K4").Value = 6003
$ws.Range("M4").Value = -5891

$ws.Range("H75").Value = 35157.75
$ws.Range("J75").Value = 35157.75
$ws.Range("L75").Value = 35157.75
$ws.Range("N75").Value = -36905.75

$ws.Range("H78").Value = 35157.75
$ws.Range("J78").Value = 35157.75
$ws.Range("L78").Value = 105473.25
$ws.Range("N78").Value = -114209.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2286.4075
$ws.Range("I16").Value = 2272.3635
$ws.Range("J16").Value = 2348.2
$ws.Range("K16").Value = 2272.3635
$ws.Range("L16").Value = 2348.2
$ws.Range("M16").Value = -2102.3635
$ws.Range("N16").Value = -2688.2

$ws.Range("H68").Value = 2318.182
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2318.182
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2318.182
$ws.Range("N68").Value = -3816.182
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 2318.182
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2318.182
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 11590.91
$ws.Range("N71").Value = -19078.91
$ws.Range("M71").ClearContents()

$ws.Range("H93").Value = 1950.037
$ws.Range("I93").Value = 1292.375
$ws.Range("J93").Value = 2226.9473
$ws.Range("K93").Value = 1292.375
$ws.Range("L93").Value = 2226.9473
$ws.Range("M93").Value = -44.375
$ws.Range("N93").Value = -4722.9473

$ws.Range("H100").Value = 2070.2856
$ws.Range("I100").Value = 1864.5555
$ws.Range("J100").Value = 2440.6
$ws.Range("K100").Value = 1864.5555
$ws.Range("L100").Value = 2440.6
$ws.Range("M100").Value = -1323.5555
$ws.Range("N100").Value = -3522.6

$ws.Range("H122").Value = 2726.25
$ws.Range("I122").Value = 2633.3333
$ws.Range("K122").Value = 7899.999899999999
$ws.Range("M122").Value = -5449.999899999999

$ws.Range("H132").Value = 2251.7258
$ws.Range("I132").Value = 1684.0513
$ws.Range("J132").Value = 3214.3044
$ws.Range("K132").Value = 5052.1539
$ws.Range("L132").Value = 9642.913199999999
$ws.Range("M132").Value = -2522.1539
$ws.Range("N132").Value = -14702.9132

$ws.Range("H136").Value = 1120.6909
$ws.Range("I136").Value = 873.4375
$ws.Range("J136").Value = 2816.1428
$ws.Range("K136").Value = 2620.3125
$ws.Range("L136").Value = 8448.428400000001
$ws.Range("M136").Value = -70.3125
$ws.Range("N136").Value = -13548.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 7300600
$ws.Range("J5").Value = 8715143
$ws.Range("L5").Value = 8715143
$ws.Range("N5").Value = -8715367

$ws.Range("H100").Value = 613.8333
$ws.Range("I100").Value = 536.8
$ws.Range("J100").Value = 999
$ws.Range("K100").Value = 1073.6
$ws.Range("L100").Value = 1998
$ws.Range("M100").Value = -532.5999999999999
$ws.Range("N100").Value = -3080

$ws.Range("H132").Value = 1848.674
$ws.Range("I132").Value = 1541.0571
$ws.Range("J132").Value = 2827.4546
$ws.Range("K132").Value = 4623.1713
$ws.Range("L132").Value = 8482.363799999999
$ws.Range("M132").Value = -2093.1713
$ws.Range("N132").Value = -13542.3638
